$wb = $excel.ActiveWorkbook

$wsTrafo = $wb.Worksheets.Item("trafo")
$wsTrafo.Range("A2:I10").ClearContents()
